$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (shifts existing D:K -> F:M)
$ws.Columns("D:E").Insert()

# Copy number formats/styles from the (now shifted) F:G columns into the
# newly inserted D:E columns, for each contiguous data block.
$ws.Range("F7:G35").Copy()
$ws.Range("D7:E35").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("F38:G77").Copy()
$ws.Range("D38:E77").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("F80:G102").Copy()
$ws.Range("D80:E102").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$excel.CutCopyMode = $false

# Populate the two newly-inserted quarter columns (D = latest quarter,
# E = second latest quarter) with their reported values.
$ws.Range("D7").Value = [DateTime]"2018-12-31"
$ws.Range("E7").Value = [DateTime]"2018-09-30"
$ws.Range("D8").Value = 69200
$ws.Range("E8").Value = 142800
$ws.Range("D9").Value = 21500
$ws.Range("E9").Value = 25400
$ws.Range("D10").Value = 47700
$ws.Range("E10").Value = 117400
$ws.Range("D12").Value = 27100
$ws.Range("E12").Value = 22900
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 118500
$ws.Range("E17").Value = 114500
$ws.Range("D18").Value = -49300
$ws.Range("E18").Value = 28300
$ws.Range("D20").Value = 1100
$ws.Range("E20").Value = 1200
$ws.Range("D21").Value = -45800
$ws.Range("E21").Value = 31900
$ws.Range("D22").Value = 5300
$ws.Range("E22").Value = 5400
$ws.Range("D23").Value = -53500
$ws.Range("E23").Value = 24100
$ws.Range("D24").Value = -48300
$ws.Range("E24").Value = 38000
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = -5200
$ws.Range("E26").Value = -13900
$ws.Range("D27").Value = -5200
$ws.Range("E27").Value = -13900
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 14800
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -1100
$ws.Range("E32").Value = -1200
$ws.Range("D33").Value = 9600
$ws.Range("E33").Value = -13900
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 9600
$ws.Range("E35").Value = -13900
$ws.Range("D38").Value = [DateTime]"2018-12-31"
$ws.Range("E38").Value = [DateTime]"2018-09-30"
$ws.Range("D41").Value = 293600
$ws.Range("E41").Value = 321000
$ws.Range("D42").Value = 152000
$ws.Range("E42").Value = 139900
$ws.Range("D43").Value = 23400
$ws.Range("E43").Value = 51500
$ws.Range("D44").Value = 29000
$ws.Range("E44").Value = 10800
$ws.Range("D45").Value = 30100
$ws.Range("E45").Value = 23400
$ws.Range("D46").Value = 528100
$ws.Range("E46").Value = 546600
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 60500
$ws.Range("E48").Value = 52100
$ws.Range("D49").Value = 710600
$ws.Range("E49").Value = 712000
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 400
$ws.Range("E52").Value = 400
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 1299700
$ws.Range("E54").Value = 1311100
$ws.Range("D57").Value = 48900
$ws.Range("E57").Value = 29600
$ws.Range("D58").Value = 600
$ws.Range("E58").Value = 600
$ws.Range("D59").Value = 90800
$ws.Range("E59").Value = 105500
$ws.Range("D60").Value = 140300
$ws.Range("E60").Value = 135700
$ws.Range("D61").Value = 343100
$ws.Range("E61").Value = 340800
$ws.Range("D62").Value = 204300
$ws.Range("E62").Value = 235900
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 687700
$ws.Range("E66").Value = 712500
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = -393800
$ws.Range("E72").Value = -403400
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 612000
$ws.Range("E76").Value = 598600
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = [DateTime]"2018-12-31"
$ws.Range("E80").Value = [DateTime]"2018-09-30"
$ws.Range("D81").Value = 9600
$ws.Range("E81").Value = -13900
$ws.Range("D83").Value = 2400
$ws.Range("E83").Value = 2500
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = -12300
$ws.Range("E89").Value = 78200
$ws.Range("D91").Value = -10800
$ws.Range("E91").Value = -11800
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -14800
$ws.Range("E94").Value = -2700
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = 100
$ws.Range("E100").Value = 2300
$ws.Range("D101").Value = -200
$ws.Range("E101").Value = -200
$ws.Range("D102").Value = -27300
$ws.Range("E102").Value = 77600
$ws.Range("H89").Value = 60500
$ws.Range("I89").Value = 51700
$ws.Range("F91").Value = -6000
$ws.Range("H91").Value = -3300
$ws.Range("I91").Value = -1600
$ws.Range("J91").Value = -3000
$ws.Range("H102").Value = 114900
$ws.Range("I102").Value = 51700
